$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("B2").Value = 79001
$ws.Range("B3").Value = 91808
$ws.Range("B4").Value = 92106
$ws.Range("B5").Value = 93095
$ws.Range("B6").Value = 79243
$ws.Range("B7").Value = 92106
$ws.Range("B8").Value = 92106
